# Generate Report for handback
#
# Marks the "dae061a2-1846-46c0-84e2-6eac767f7636.md" file as handed back
# (for both the zh-cn and de-de locales), updating the Status columns and
# recording new "Latest Handback DateTime" values.

$wb = $excel.ActiveWorkbook

$statusHandedBack = "Handed back: in sync with en-US"

# --- Overview sheet: mirror the status for both locales ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = $statusHandedBack
$wsOverview.Range("C3").Value = $statusHandedBack

# --- zh-cn sheet: update status + new handback datetime for row 3 ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("B3").Value = $statusHandedBack
$wsZhCn.Range("G3").Value = "2016-01-25 13:43:03"

# --- de-de sheet: update status + new handback datetime for row 3 ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("B3").Value = $statusHandedBack
$wsDeDe.Range("G3").Value = "2016-01-25 13:43:24"
